$d = $word.ActiveDocument

# Remove the first paragraph "DDoS + AI = NAPATZ" and the following empty paragraph.
$p1 = $d.Paragraphs(1)
$p1.Range.Delete()

$p2 = $d.Paragraphs(1)
$p2.Range.Delete()
